$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "Field width specifications in printf" paragraphs by content,
# so the script is resilient to any paragraph-numbering differences.
# ---------------------------------------------------------------------------
$pLeadingIdx = 0
$pCodeIdx = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($pLeadingIdx -eq 0 -and $t -match "^If the integer only takes one column") {
        $pLeadingIdx = $i
    }
    if ($t -match "printf\(.%-2i\\n.\, n\);") {
        $pCodeIdx = $i
    }
}

# ---------------------------------------------------------------------------
# 1) "If the integer only takes one column..." paragraph: the sentence was
#    originally split into two runs around a stray _GoBack bookmark. Replace
#    the whole paragraph's text with the complete sentence as a single run,
#    which also removes the now-stale bookmark from this spot.
# ---------------------------------------------------------------------------
$pLeading = $d.Paragraphs($pLeadingIdx)
$rLeading = $pLeading.Range
$leadingStart = $rLeading.Start
$leadingTextEnd = $rLeading.End - 1   # exclude the paragraph mark

$d.Range($leadingStart, $leadingTextEnd).Text = ""
$d.Range($leadingStart, $leadingStart).Text = "If the integer only takes one column, then a leading space will be displayed along with the integer in order to fill the two columns. This is similar for other field width specifications."

# ---------------------------------------------------------------------------
# 2) Right after the "printf(“%-2i\n”, n);" left-justified example, add two
#    new paragraphs: an explanatory sentence (carrying the relocated
#    _GoBack bookmark at the point the author's cursor last was) and a new
#    code sample demonstrating the "%.3i" field width specifier.
# ---------------------------------------------------------------------------
$pCode = $d.Paragraphs($pCodeIdx)
$rCode = $pCode.Range
$rCode.InsertParagraphAfter()

$pExplain = $d.Paragraphs($pCodeIdx + 1)
$explainXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:r><w:t xml:space='preserve'>You can also fill the </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>empty </w:t></w:r>" +
    "<w:r><w:t>colum</w:t></w:r>" +
    "<w:r><w:t>n</w:t></w:r>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
    "<w:r><w:t>s with zeros using a dot character, for right justification:</w:t></w:r>" +
    "</w:p>"
$null = $pExplain.Range.InsertXML($explainXml)

$pExplain = $d.Paragraphs($pCodeIdx + 1)
$rExplain = $pExplain.Range
$rExplain.InsertParagraphAfter()

$pNewCode = $d.Paragraphs($pCodeIdx + 2)
$codeXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
    "<w:pPr><w:pStyle w:val='code'/></w:pPr>" +
    "<w:r><w:t>printf(&#8220;%.3i\n&#8221;, n);</w:t></w:r>" +
    "</w:p>"
$null = $pNewCode.Range.InsertXML($codeXml)

Write-Output "Applied field-width-specifier edits."
